$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("D2").Value2 = "'20220303"
Write-Host "done"
